# Weekly update: a new daily price record is inserted at row 77
# (Provincia de Linares, fecha 45079), shifting all subsequent rows
# down by one. The former last row becomes row 171.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row before row 77; existing rows 77:170 shift to 78:171.
$ws.Rows("77").Insert()

# Populate the newly inserted row with the new record's data.
$ws.Range("A77").Value = 9
$ws.Range("B77").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C77").Value = "Metropolitana"
$ws.Range("D77").Value = 45079
$ws.Range("E77").Value = 13
$ws.Range("F77").Value = "Fruta"
$ws.Range("G77").Value = 100101
$ws.Range("H77").Value = "Berries"
$ws.Range("I77").Value = 100101004
$ws.Range("J77").Value = "Frambuesa"
$ws.Range("K77").Value = "Sin especificar"
$ws.Range("L77").Value = "Primera"
$ws.Range("M77").Value = 500
$ws.Range("N77").Value = 9000
$ws.Range("O77").Value = 9500
$ws.Range("P77").Value = 9280
$ws.Range("Q77").Value = "$/bandeja 2 kilos"
$ws.Range("R77").Value = "Provincia de Linares"
$ws.Range("S77").Value = 4640
$ws.Range("T77").Value = 2
